# "added first version of DESeq2 analysis in pipeline"
#
# The DESeq2 section of the sheet (rows 30-33) already references the
# "featureCountMethod" example value in B5 and the "countTableOrigin"
# example value in B31. Both used to read the placeholder "featureCounts";
# this edit switches both examples over to "HTSeq" instead, and leaves the
# cursor/selection on the last cell of the DESeq2 block (B33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "HTSeq"
$ws.Range("B31").Value = "HTSeq"

$ws.Range("B33").Select()
